# Applies the LOM3243.xlsx restructuring described in the commit diff:
#  - sheet shrinks from 25 data rows to 23 (dimension A1:C25 -> A1:C23)
#  - a "Programa resumido:/Short syllabus:" block is inserted after the
#    "Docentes responsaveis:" rows, shifting everything below it down,
#    while the trailing "Requisitos:" detail row is removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the last two rows first so the grid is the correct final size (23 rows)
# before we rewrite the content; deleting whole rows (rather than just clearing
# cell contents) is what shrinks <dimension> from C25 down to C23.
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(24).Delete()

# Row 10
$ws.Range("B10").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C10").Value = "3577649 - Carlos Angelo Nunes"

# Row 13
$ws.Range("A13").Value = "Programa resumido:"
# Assigning the literal text "01/01/2012" via .Value would be auto-parsed as a
# real date by Excel's type coercion; copy/paste the existing text cell (B8/C8,
# which already stores this same shared string as plain text) instead so the
# destination keeps its text type (and its own column style, B13/C13 untouched).
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial()
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial()
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range("A16").Value = "Syllabus:"
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).AutoFit()  # drop custom height back to default

# Row 18
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Os seminários proferidos por estudantes de graduação e pós-graduação, professores e convidados serão debatidos e analisados pelos alunos em forma de relatório. Os seminários apresentados pelos alunos serão avaliados na disciplina."
$ws.Range("C19").Value = "Os seminários proferidos por estudantes de graduação e pós-graduação, professores e convidados serão debatidos e analisados pelos alunos em forma de relatório. Os seminários apresentados pelos alunos serão avaliados na disciplina."
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A nota final será calculada pela média aritmética dos relatórios e do seminário."
$ws.Range("C20").Value = "A nota final será calculada pela média aritmética dos relatórios e do seminário."

# Row 21
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "Não há."
$ws.Range("C21").Value = "Não há."
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Rows.Item(22).AutoFit()  # drop custom height back to default

# Row 23
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOM3238 -  Projeto Integrado I  (Requisito)`n"
$ws.Range("C23").Value = "LOM3238 -  Projeto Integrado I  (Requisito)`n"
$ws.Rows.Item(23).RowHeight = 30

